# Weekly refresh of the "Hortaliza, Macroferia Regional de Talca - Betarraga"
# price series: a new week's record is inserted at the top of the data block
# (row 375, right after the header + first data rows), pushing the existing
# rows 375-414 down to 376-415 and growing the used range from A1:R414 to
# A1:R415.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 375 (shifts 375:414 -> 376:415).
$ws.Range("A375").EntireRow.Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A375").Value = 5
$ws.Range("B375").Value = "Macroferia Regional de Talca"
$ws.Range("C375").Value = "Maule"
$ws.Range("D375").Value = 44858
$ws.Range("E375").Value = 7
$ws.Range("F375").Value = 100114014
$ws.Range("G375").Value = "Betarraga"
$ws.Range("H375").Value = "Sin especificar"
$ws.Range("I375").Value = "Segunda"
$ws.Range("J375").Value = 4000
$ws.Range("K375").Value = 900
$ws.Range("L375").Value = 900
$ws.Range("M375").Value = 900
$ws.Range("N375").Value = "`$/paquete 5 unidades"
$ws.Range("O375").Value = "Región del Maule"
$ws.Range("P375").Value = 180
$ws.Range("Q375").Value = 5
$ws.Range("R375").Value = "Hortaliza"
